# Update the "datos actualizados" timestamp in the title row (A1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 13 de Agosto de 2020 a las 09:52"

# --- Update per-country statistics (columns: B=Casos totales, C=Nuevos casos,
#     D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes) ---

# Row 6: India
$ws.Cells.Item(6, 2).Value = 2399992
$ws.Cells.Item(6, 3).Value = 4521
$ws.Cells.Item(6, 4).Value = 1697811
$ws.Cells.Item(6, 5).Value = 655012
$ws.Cells.Item(6, 7).Value = 31
$ws.Cells.Item(6, 8).Value = 47169

# Row 7: Rusia
$ws.Cells.Item(7, 2).Value = 907758
$ws.Cells.Item(7, 3).Value = 5057
$ws.Cells.Item(7, 4).Value = 716396
$ws.Cells.Item(7, 5).Value = 175978
$ws.Cells.Item(7, 7).Value = 124
$ws.Cells.Item(7, 8).Value = 15384

# Row 47: Singapur
$ws.Cells.Item(47, 2).Value = 55497
$ws.Cells.Item(47, 3).Value = 102
$ws.Cells.Item(47, 5).Value = 4950

# Row 56: Armenia
$ws.Cells.Item(56, 2).Value = 41023
$ws.Cells.Item(56, 3).Value = 229
$ws.Cells.Item(56, 4).Value = 33897
$ws.Cells.Item(56, 5).Value = 6317
$ws.Cells.Item(56, 7).Value = 3
$ws.Cells.Item(56, 8).Value = 809

# Row 108: Hungria
$ws.Cells.Item(108, 2).Value = 4813
$ws.Cells.Item(108, 3).Value = 45
$ws.Cells.Item(108, 4).Value = 3561
$ws.Cells.Item(108, 5).Value = 645
$ws.Cells.Item(108, 7).Value = 2
$ws.Cells.Item(108, 8).Value = 607

# Row 122: Sri Lanka
$ws.Cells.Item(122, 4).Value = 2646
$ws.Cells.Item(122, 5).Value = 224

# Row 201: Fiyi
$ws.Cells.Item(201, 4).Value = 20
$ws.Cells.Item(201, 5).Value = 7

# Rows 213/214: re-sort "Montserrat" / "Islas Malvinas" alphabetically
# (Islas Malvinas now comes before Montserrat), swapping the full row
# contents between the two rows.
$row213 = @($ws.Cells.Item(213, 1).Value2, $ws.Cells.Item(213, 2).Value2, $ws.Cells.Item(213, 3).Value2, $ws.Cells.Item(213, 4).Value2, $ws.Cells.Item(213, 5).Value2, $ws.Cells.Item(213, 6).Value2, $ws.Cells.Item(213, 7).Value2, $ws.Cells.Item(213, 8).Value2)
$row214 = @($ws.Cells.Item(214, 1).Value2, $ws.Cells.Item(214, 2).Value2, $ws.Cells.Item(214, 3).Value2, $ws.Cells.Item(214, 4).Value2, $ws.Cells.Item(214, 5).Value2, $ws.Cells.Item(214, 6).Value2, $ws.Cells.Item(214, 7).Value2, $ws.Cells.Item(214, 8).Value2)

for ($c = 1; $c -le 8; $c++) {
    $ws.Cells.Item(213, $c).Value = $row214[$c - 1]
    $ws.Cells.Item(214, $c).Value = $row213[$c - 1]
}
